# The deck's Design/theme (ppt/theme/theme1.xml, "Integral" / "Red Violet")
# is switched to the built-in default "Office" color palette - the same
# operation as opening the Design tab and clicking the "Office" theme
# color swatch. PowerPoint exposes the 12 theme colors (dk1, lt1, dk2,
# lt2, accent1-6, hlink, folHlink - in that fixed order) through
# Design.SlideMaster.Theme.ThemeColorScheme.Item(n).RGB, each value
# being a standard OLE RGB() long (R + G*256 + B*65536).

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

# Office theme colors, in Item() order: dk1, lt1, dk2, lt2,
# accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$officeRgb = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeRgb[$i - 1]
}
